$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.222.89"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "2.959.32"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'379.60"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").Value = "'102.62"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("D7").Value = "'0.541"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.589"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").Value = "'36.55"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "'0.0839"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "3.422.14"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").Value = "'18.01"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").Value = "'7.38"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "2.940.07"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "'0.985"
$ws.Range("E17").Value = "  +4.95%  "
$ws.Range("D18").Value = "51.241.78"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "'3.22"
$ws.Range("E19").Value = "  -5.88%  "
$ws.Range("D20").Value = "'7.18"
$ws.Range("E20").Value = "  -2.05%  "
$ws.Range("D21").Value = "'12.58"
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("D22").Value = "0.0₃0952"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "'68.35"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'261.67"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'2.86"
$ws.Range("E25").Value = "  +2.36%  "
$ws.Range("D26").Value = "'8.28"
$ws.Range("E26").Value = "  +13.06%  "
$ws.Range("D27").Value = "'7.60"
$ws.Range("E27").Value = "  +7.65%  "
$ws.Range("D28").Value = "'0.169"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "'25.72"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("E32").Value = "  +10.12%  "
$ws.Range("D33").Value = "'9.79"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").Value = "'50.49"
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'33.78"
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("B36").Value = "Toncoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D36").Value = "'2.04"
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("D37").Value = "'0.0444"
$ws.Range("E37").Value = "  +4.18%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'2.98"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("D40").Value = "'16.92"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").Value = "'2.56"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").Value = "'1.78"
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("D44").Value = "'122.07"
$ws.Range("E44").Value = "  -1.86%  "
$ws.Range("D45").Value = "'21.01"
$ws.Range("E45").Value = "  -4.13%  "
$ws.Range("D46").Value = "'2.05"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("D47").Value = "'0.272"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("D49").Value = "2.007.92"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").Value = "'3.21"
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("D51").Value = "'0.0335"
$ws.Range("E51").Value = "  +4.47%  "
